$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text change: "lb" -> "Lower bound" ---
$ws.Range("B1").Value = "Lower bound"

# --- Column widths: A gets a width, B gets a new (wider) width ---
$ws.Columns.Item(1).ColumnWidth = 8.6640625
$ws.Columns.Item(2).ColumnWidth = 11.83203125

# --- Selection moves to H23 ---
$ws.Range("H23").Select()

# --- Window size grows ---
$excel.Width = 1037
$excel.Height = 543

# --- Formatting: Times New Roman font, centered, for the whole table ---
$fullRange = $ws.Range("A1:E5")
$fullRange.Font.Name = "Times New Roman"
$fullRange.HorizontalAlignment = -4108

# --- Header row gets a boxed (top+bottom) thin border ---
$hdrRow = $ws.Range("A1:E1")
$hdrRow.Borders.Item(9).LineStyle = 1
$hdrRow.Borders.Item(8).LineStyle = 1

# --- Last data row gets a bottom thin border (closes the table) ---
$lastRow = $ws.Range("A5:E5")
$lastRow.Borders.Item(9).LineStyle = 1
